# Update counter values (column F) on sheets "展览" (sheet1) and "全部类型" (sheet4)
# to reflect newly generated output, per commit "Update gh-pages to output
# generated at 456a3b4".

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsAllTypes = $wb.Worksheets.Item("全部类型")

# Sheet "展览" (sheet1.xml) updates
$wsExhibit.Range("F4").Value = 1280
$wsExhibit.Range("F7").Value = 985
$wsExhibit.Range("F12").Value = 299
$wsExhibit.Range("F16").Value = 4055
$wsExhibit.Range("F19").Value = 2646
$wsExhibit.Range("F21").Value = 1086
$wsExhibit.Range("F23").Value = 775
$wsExhibit.Range("F25").Value = 43
$wsExhibit.Range("F26").Value = 2331
$wsExhibit.Range("F30").Value = 689
$wsExhibit.Range("F40").Value = 286

# Sheet "全部类型" (sheet4.xml) updates
$wsAllTypes.Range("F4").Value = 1280
$wsAllTypes.Range("F6").Value = 985
$wsAllTypes.Range("F16").Value = 4055
$wsAllTypes.Range("F20").Value = 2646
$wsAllTypes.Range("F21").Value = 1086
$wsAllTypes.Range("F23").Value = 775
$wsAllTypes.Range("F26").Value = 43
$wsAllTypes.Range("F27").Value = 2331
$wsAllTypes.Range("F35").Value = 689
$wsAllTypes.Range("F45").Value = 286

$wb.Save()
